$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cam_FilingPeriodLobbyist")

$row = 31
$srcRow = 30

# Reuse the existing date-number-format style (already applied to columns
# B, C and G on the prior row) instead of creating new style entries.
foreach ($col in 2, 3, 7) {
    $ws.Cells.Item($srcRow, $col).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 30
$ws.Cells.Item($row, 2).Value = 43480.4993055556
$ws.Cells.Item($row, 3).Value = 43465.4993055556
$ws.Cells.Item($row, 4).Value = "2019 January Lobbyist Report"
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 43375

# Extend the workbook-level defined name to cover the new row.
$name = $wb.Names.Item("Cam_FilingPeriodLobbyist")
$name.RefersTo = "=Cam_FilingPeriodLobbyist!`$A`$1:`$H`$31"
